$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-run docking results: tiny refinements to a handful of raw values
# (re-running the docking pass nudges the last bit of a few numbers) ---
$ws.Range("C10").Value = 157.70599999999996
$ws.Range("D10").Value = 146.389999999999958
$ws.Range("D18").Value = 157.877999999999957
$ws.Range("E22").Value = 153.776999999999958
$ws.Range("E25").Value = 154.44999999999996
$ws.Range("D31").Value = 159.69999999999996

# --- Normalize the duplicate "header row" styles (C1:E1, C15:E15, C30:E30)
# so they all collapse back onto the same shared cell style ---
$headerRanges = @("C1:E1", "C15:E15", "C30:E30")
foreach ($rng in $headerRanges) {
    $ws.Range($rng).Font.Name = "Arial"
    $ws.Range($rng).Font.Bold = $true
    $ws.Range($rng).Font.Size = 10
    $ws.Range($rng).HorizontalAlignment = -4108
}

# --- Normalize the duplicate "row label" style used by B16:B27 so it
# collapses back onto the same shared style as B2:B13 / B31:B36 ---
$ws.Range("B16:B27").Font.Name = "Arial"
$ws.Range("B16:B27").Font.Bold = $true
$ws.Range("B16:B27").Font.Size = 10

# --- Widen the X / Y data columns for readability ---
$ws.Range("C:C").ColumnWidth = 15.333333333333334
$ws.Range("D:D").ColumnWidth = 15.0

# --- Scroll the view down and move the selection ---
$ws.Range("D25").Select()
$ws.Application.ActiveWindow.ScrollRow = 17
